## Commit: Fri, Jul 31, 2020  5:05:24 PM
##
## Re-style the three data tables (slides 14, 15 and 16) with a different
## built-in table style, as done from the Table Design ribbon
## ("Table Styles" gallery) in PowerPoint.

$p = $ppt.ActivePresentation

$newTableStyleId = "{528E511D-A9D6-42A0-AD7D-873BA1EBE824}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)

    # The table on each of these slides is the first shape (a GraphicFrame
    # hosting a Table).
    $shape = $slide.Shapes.Item(1)

    if ($shape.HasTable) {
        $table = $shape.Table
        # Table styles can't be assigned through the Style property directly;
        # PowerPoint exposes this as Table.ApplyStyle(styleId).
        $table.ApplyStyle($newTableStyleId)
    }
}
